# Auto-generated edit script applying the diff to 上海-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3547
$ws1.Range("F5").Value = 8296
$ws1.Range("F13").Value = 63
$ws1.Range("F14").Value = 629
$ws1.Range("F16").Value = 7265
$ws1.Range("F17").Value = 451
$ws1.Range("F18").Value = 7558
$ws1.Range("F20").Value = 57260
$ws1.Range("F21").Value = 57261
$ws1.Range("F22").Value = 4667
$ws1.Range("F26").Value = 480
$ws1.Range("F33").Value = 94
$ws1.Range("F35").Value = 887
$ws1.Range("F36").Value = 1300
$ws1.Range("F37").Value = 1629
$ws1.Range("F39").Value = 174
$ws1.Range("F41").Value = 1082
$ws1.Range("F47").Value = 189

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 34
$ws2.Range("F8").Value = 48
$ws2.Range("F28").Value = 2
$ws2.Range("F33").Value = 2
$ws2.Range("F34").Value = 1
$ws2.Range("F40").Value = 113
$ws2.Range("F41").Value = 197
$ws2.Range("F47").Value = 271

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2344
$ws3.Range("F7").Value = 675
$ws3.Range("F8").Value = 2398
$ws3.Range("F12").Value = 105
$ws3.Range("F15").Value = 246
$ws3.Range("F16").Value = 2191
$ws3.Range("F17").Value = 35
$ws3.Range("F18").Value = 461
$ws3.Range("G8").Value = "已售罄"

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value = "2024-09-06"
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value = "2024-09-10"
$ws4.Range("C6").Value = "上海·「HUNTER×HUNTER × animate cafe」"
$ws4.Range("C7").Value = "上海·迷你四驱车赛场"
$ws4.Range("D6").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws4.Range("D7").Value = "虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶"
$ws4.Range("E6").Value = "2024.09.06 00:00-10.08 23:59"
$ws4.Range("E7").Value = "2024.09.10 10:00-12.31 22:00"
$ws4.Range("F2").Value = 2344
$ws4.Range("F4").Value = 8296
$ws4.Range("F5").Value = 675
$ws4.Range("F6").Value = 1731
$ws4.Range("F7").Value = 4
$ws4.Range("F10").Value = 2191
$ws4.Range("F11").Value = 63
$ws4.Range("F12").Value = 7558
$ws4.Range("F13").Value = 57262
$ws4.Range("F15").Value = 34
$ws4.Range("F17").Value = 4667
$ws4.Range("F24").Value = 4972
$ws4.Range("F25").Value = 48
$ws4.Range("F26").Value = 94
$ws4.Range("F28").Value = 887
$ws4.Range("F29").Value = 1300
$ws4.Range("F30").Value = 1629
$ws4.Range("F32").Value = 461
$ws4.Range("F36").Value = 174
$ws4.Range("F38").Value = 1082
$ws4.Range("F45").Value = 2
$ws4.Range("F51").Value = 271
$ws4.Range("G6").Value = 30
$ws4.Range("G7").Value = 48
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91069"
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=92042"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202408/4GkLI2cn1724227065219.jpeg"
$ws4.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202409/LzFT5TMO1725348229429.png"

